# Generate Report for Handback
# Update handback-status workbook with new file identifiers / timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$wsOverview.Range("A2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"
$wsOverview.Range("B2").Value = "e2e\714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"
$wsOverview.Range("G2").Value = "2016-09-03 07:06:31"

$wsOverview.Range("A3").Value = "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md"
$wsOverview.Range("B3").Value = "e2e\ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md"
$wsOverview.Range("G3").Value = "2016-09-03 07:06:31"

# ---- zh-cn sheet ----
$wsZhCn.Range("A2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"
$wsZhCn.Range("G2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-03 07:06:26"
$wsZhCn.Range("I2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"
$wsZhCn.Range("J2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-03 07:06:43"

$wsZhCn.Range("A3").Value = "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md"
$wsZhCn.Range("G3").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-03 07:06:26"
$wsZhCn.Range("I3").Value = "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md"
$wsZhCn.Range("J3").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-03 07:06:43"

# ---- de-de sheet ----
$wsDeDe.Range("A2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"
$wsDeDe.Range("G2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-03 07:06:31"
$wsDeDe.Range("I2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md"
$wsDeDe.Range("J2").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-03 07:06:53"

$wsDeDe.Range("A3").Value = "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md"
$wsDeDe.Range("G3").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-03 07:06:31"
$wsDeDe.Range("I3").Value = "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md"
$wsDeDe.Range("J3").Value = "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.0bb80c2d91da5d941be62748220d706d7ed76718.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-03 07:06:53"

# ---- Keep hyperlink display text (and underlying file names) in sync ----
# In-place hyperlink edits aren't supported by this runtime (they would
# duplicate the <hyperlink> entry), so rebuild each sheet's hyperlinks:
# drop the old ones and re-add them pointing at the same target URLs with
# the refreshed display text.

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990020c3aa22be0c42647eafc66485ac2263fd1a/e2e/714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md", "", "", "e2e\714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990020c3aa22be0c42647eafc66485ac2263fd1a/e2e/ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md", "", "", "e2e\ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md") | Out-Null

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990020c3aa22be0c42647eafc66485ac2263fd1a/e2e/714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md", "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8b15e9acd6cf8aff81ee5feb2bdc8d1746fe7998/e2e/714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md", "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990020c3aa22be0c42647eafc66485ac2263fd1a/e2e/ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md", "", "", "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8b15e9acd6cf8aff81ee5feb2bdc8d1746fe7998/e2e/ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md", "", "", "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md") | Out-Null

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990020c3aa22be0c42647eafc66485ac2263fd1a/e2e/714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md", "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cb59eca4446dc5ccaf44fc115a691c588ccf1896/e2e/714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md", "", "", "714c55c5-6b7d-4e52-b6e7-6f7268cee2a7.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990020c3aa22be0c42647eafc66485ac2263fd1a/e2e/ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md", "", "", "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cb59eca4446dc5ccaf44fc115a691c588ccf1896/e2e/ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md", "", "", "ffff2a6f45f1-303c-4443-8918-86936b0b4b4e.md") | Out-Null
